$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 831.83636
$ws.Range("I15").Value = 831.83636
$ws.Range("K15").Value = 2495.50908
$ws.Range("M15").Value = -2326.50908
$ws.Range("H51").Value = 4748
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 5372.6665
$ws.Range("K51").Value = 1000
$ws.Range("L51").Value = 5372.6665
$ws.Range("M51").Value = -516
$ws.Range("N51").Value = -6340.6665
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()
$ws.Range("H61").Value = 408.5
$ws.Range("I61").Value = 408.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1225.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1053.5
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 2934.6667
$ws.Range("I62").Value = 2934.6667
$ws.Range("K62").Value = 2934.6667
$ws.Range("M62").Value = -2310.6667
$ws.Range("H65").Value = 2934.6667
$ws.Range("I65").Value = 2934.6667
$ws.Range("K65").Value = 14673.3335
$ws.Range("M65").Value = -11553.3335
$ws.Range("H69").Value = 7083.3887
$ws.Range("J69").Value = 8072.143
$ws.Range("L69").Value = 24216.429
$ws.Range("N69").Value = -25964.429
$ws.Range("H70").Value = 9168.625
$ws.Range("J70").Value = 9478.143
$ws.Range("L70").Value = 28434.429
$ws.Range("N70").Value = -28974.429
$ws.Range("H72").Value = 7083.3887
$ws.Range("J72").Value = 8072.143
$ws.Range("L72").Value = 72649.287
$ws.Range("N72").Value = -81385.287
$ws.Range("H73").Value = 9168.625
$ws.Range("J73").Value = 9478.143
$ws.Range("L73").Value = 28434.429
$ws.Range("N73").Value = -30306.429
$ws.Range("H106").Value = 2999.3
$ws.Range("I106").Value = 2943.6667
$ws.Range("K106").Value = 2943.6667
$ws.Range("M106").Value = -2312.6667
$ws.Range("H111").Value = 750
$ws.Range("I111").Value = 750
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2250
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 817
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 6326.4
$ws.Range("I113").Value = 2996.3333
$ws.Range("J113").Value = 11321.5
$ws.Range("K113").Value = 2996.3333
$ws.Range("L113").Value = 11321.5
$ws.Range("M113").Value = 257.6667000000002
$ws.Range("N113").Value = -17829.5
$ws.Range("H123").Value = 82000
$ws.Range("J123").Value = 82000
$ws.Range("L123").Value = 82000
$ws.Range("N123").Value = -91800
$ws.Range("H125").Value = 6188.8125
$ws.Range("I125").Value = 7444.6665
$ws.Range("J125").Value = 2421.25
$ws.Range("K125").Value = 67001.9985
$ws.Range("L125").Value = 21791.25
$ws.Range("M125").Value = -64541.9985
$ws.Range("N125").Value = -26711.25
$ws.Range("H132").Value = 834947.5600000001
$ws.Range("I132").Value = 1836.5
$ws.Range("K132").Value = 5509.5
$ws.Range("M132").Value = -2979.5
$ws.Range("H137").Value = 4423.8
$ws.Range("I137").Value = 2811.3635
$ws.Range("J137").Value = 6394.5557
$ws.Range("K137").Value = 8434.0905
$ws.Range("L137").Value = 19183.6671
$ws.Range("M137").Value = -5884.0905
$ws.Range("N137").Value = -24283.6671
$ws.Range("H138").Value = 8128.4443
$ws.Range("I138").Value = 4497.2856
$ws.Range("J138").Value = 9399.35
$ws.Range("K138").Value = 13491.8568
$ws.Range("L138").Value = 28198.05
$ws.Range("M138").Value = -8351.856800000001
$ws.Range("N138").Value = -38478.05
$ws.Range("H141").Value = 3923.8333
$ws.Range("I141").Value = 3703.1738
$ws.Range("K141").Value = 11109.5214
$ws.Range("M141").Value = -5929.5214

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8804.492
$ws.Range("I32").Value = 3140.9805
$ws.Range("K32").Value = 3140.9805
$ws.Range("M32").Value = -2853.9805
$ws.Range("H45").Value = 3749.75
$ws.Range("I45").Value = 3749.75
$ws.Range("K45").Value = 3749.75
$ws.Range("M45").Value = -3372.75
$ws.Range("H61").Value = 2684.862
$ws.Range("I61").Value = 2130.8696
$ws.Range("K61").Value = 2130.8696
$ws.Range("M61").Value = -1918.8696
$ws.Range("H74").Value = 2913
$ws.Range("I74").Value = 1315.1666
$ws.Range("K74").Value = 1315.1666
$ws.Range("M74").Value = -441.1666
$ws.Range("H77").Value = 2913
$ws.Range("I77").Value = 1315.1666
$ws.Range("K77").Value = 6575.833000000001
$ws.Range("M77").Value = -2207.833000000001
$ws.Range("H102").Value = 2584.7856
$ws.Range("I102").Value = 1471.6364
$ws.Range("J102").Value = 6666.3335
$ws.Range("K102").Value = 1471.6364
$ws.Range("L102").Value = 6666.3335
$ws.Range("M102").Value = 150.3635999999999
$ws.Range("N102").Value = -9910.333500000001
$ws.Range("H122").Value = 2291.6128
$ws.Range("I122").Value = 2173.04
$ws.Range("K122").Value = 6519.12
$ws.Range("M122").Value = -4069.12
$ws.Range("H132").Value = 2987.3809
$ws.Range("I132").Value = 3203.6843
$ws.Range("J132").Value = 932.5
$ws.Range("K132").Value = 9611.052899999999
$ws.Range("L132").Value = 2797.5
$ws.Range("M132").Value = -7081.052899999999
$ws.Range("N132").Value = -7857.5
$ws.Range("H136").Value = 2684.862
$ws.Range("I136").Value = 2130.8696
$ws.Range("K136").Value = 6392.6088
$ws.Range("M136").Value = -3842.6088

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H94").Value = 2469
$ws.Range("I94").Value = 2429.7778
$ws.Range("J94").Value = 2586.6667
$ws.Range("K94").Value = 2429.7778
$ws.Range("L94").Value = 2586.6667
$ws.Range("M94").Value = -1978.7778
$ws.Range("N94").Value = -3488.6667
$ws.Range("H107").Value = 1729.2
$ws.Range("I107").Value = 495
$ws.Range("K107").Value = 495
$ws.Range("M107").Value = 1425
$ws.Range("H134").Value = 2410.4375
$ws.Range("I134").Value = 2197.5925
$ws.Range("J134").Value = 3559.8
$ws.Range("K134").Value = 6592.7775
$ws.Range("L134").Value = 10679.4
$ws.Range("M134").Value = -4057.7775
$ws.Range("N134").Value = -15749.4
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 67500
$ws.Range("J138").Value = 67500
$ws.Range("L138").Value = 67500
$ws.Range("N138").Value = -77780
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1148.25
$ws.Range("I7").Value = 1398
$ws.Range("K7").Value = 1398
$ws.Range("M7").Value = -1285
$ws.Range("H16").Value = 1381.25
$ws.Range("I16").Value = 1225
$ws.Range("J16").Value = 1850
$ws.Range("K16").Value = 1225
$ws.Range("L16").Value = 1850
$ws.Range("M16").Value = -938
$ws.Range("N16").Value = -2424
$ws.Range("H31").Value = 7221.049
$ws.Range("I31").Value = 5276.857
$ws.Range("K31").Value = 5276.857
$ws.Range("M31").Value = -4981.857
$ws.Range("H34").Value = 7221.049
$ws.Range("I34").Value = 5276.857
$ws.Range("K34").Value = 5276.857
$ws.Range("M34").Value = -5074.857
$ws.Range("H58").Value = 1442.081
$ws.Range("I58").Value = 1442.081
$ws.Range("K58").Value = 1442.081
$ws.Range("M58").Value = -1239.081
$ws.Range("H113").Value = 1381.25
$ws.Range("I113").Value = 1225
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 1225
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 945
$ws.Range("N113").Value = -6190
$ws.Range("H122").Value = 4932.4136
$ws.Range("I122").Value = 3930.9092
$ws.Range("J122").Value = 5544.4443
$ws.Range("K122").Value = 11792.7276
$ws.Range("L122").Value = 16633.3329
$ws.Range("M122").Value = -9342.7276
$ws.Range("N122").Value = -21533.3329
$ws.Range("H132").Value = 1514.5834
$ws.Range("I132").Value = 1248.1
$ws.Range("K132").Value = 3744.3
$ws.Range("M132").Value = -1214.3
$ws.Range("H134").Value = 2634.25
$ws.Range("I134").Value = 2531.3914
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 7594.174199999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5059.174199999999
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 1442.081
$ws.Range("I136").Value = 1442.081
$ws.Range("K136").Value = 4326.242999999999
$ws.Range("M136").Value = -1776.242999999999
$ws.Range("H137").Value = 109999
$ws.Range("J137").Value = 109999
$ws.Range("L137").Value = 109999
$ws.Range("N137").Value = -120199

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 294151.34
$ws.Range("I2").Value = 454576
$ws.Range("J2").Value = 39.5
$ws.Range("K2").Value = 2727456
$ws.Range("L2").Value = 237
$ws.Range("M2").Value = -2727343
$ws.Range("N2").Value = -463
$ws.Range("H7").Value = 719.3333
$ws.Range("I7").Value = 979.5
$ws.Range("K7").Value = 2938.5
$ws.Range("M7").Value = -2826.5
$ws.Range("H38").Value = 107.2
$ws.Range("I38").Value = 107.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 321.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 25.39999999999998
$ws.Range("N38").ClearContents()
$ws.Range("H55").Value = 462
$ws.Range("I55").Value = 462
$ws.Range("K55").Value = 1386
$ws.Range("M55").Value = -1209
$ws.Range("H104").Value = 343461.53
$ws.Range("I104").Value = 5000
$ws.Range("J104").Value = 371666.66
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 1114999.98
$ws.Range("M104").Value = -12379
$ws.Range("N104").Value = -1120241.98
$ws.Range("H129").Value = 6375.4165
$ws.Range("I129").Value = 979.5
$ws.Range("J129").Value = 9073.375
$ws.Range("K129").Value = 2938.5
$ws.Range("L129").Value = 27220.125
$ws.Range("M129").Value = 2061.5
$ws.Range("N129").Value = -37220.125
$ws.Range("H131").Value = 7038729
$ws.Range("I131").Value = 139822.25
$ws.Range("K131").Value = 419466.75
$ws.Range("M131").Value = -414426.75
$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 8991
$ws.Range("M132").Value = -6461
$ws.Range("H137").Value = 2761.2307
$ws.Range("I137").Value = 2736.8572
$ws.Range("J137").Value = 2789.6667
$ws.Range("K137").Value = 8210.571599999999
$ws.Range("L137").Value = 8369.000100000001
$ws.Range("M137").Value = -3110.571599999999
$ws.Range("N137").Value = -18569.0001
$ws.Range("H140").Value = 2048.5
$ws.Range("I140").Value = 2095.111
$ws.Range("K140").Value = 6285.333
$ws.Range("M140").Value = -1105.333
$ws.Range("H141").Value = 1411.2
$ws.Range("I141").Value = 1411.2
$ws.Range("K141").Value = 4233.6
$ws.Range("M141").Value = 946.3999999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 114.666664
$ws.Range("I2").Value = 114.666664
$ws.Range("K2").Value = 114.666664
$ws.Range("M2").Value = -1.666663999999997
$ws.Range("H122").Value = 1688.625
$ws.Range("I122").Value = 1644.1428
$ws.Range("K122").Value = 4932.428400000001
$ws.Range("M122").Value = -2482.428400000001
$ws.Range("H132").Value = 3311.842
$ws.Range("I132").Value = 3371.0625
$ws.Range("J132").Value = 2996
$ws.Range("K132").Value = 10113.1875
$ws.Range("L132").Value = 8988
$ws.Range("M132").Value = -7583.1875
$ws.Range("N132").Value = -14048
$ws.Range("H134").Value = 58729.8
$ws.Range("J134").Value = 58729.8
$ws.Range("L134").Value = 176189.4
$ws.Range("N134").Value = -181259.4
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2943.6667
$ws.Range("J7").Value = 2482.1667
$ws.Range("L7").Value = 2482.1667
$ws.Range("N7").Value = -2706.1667
$ws.Range("H10").Value = 12834.333
$ws.Range("J10").Value = 12834.333
$ws.Range("L10").Value = 12834.333
$ws.Range("N10").Value = -13114.333
$ws.Range("H22").Value = 1715.2222
$ws.Range("I22").Value = 1223.6316
$ws.Range("J22").Value = 2882.75
$ws.Range("K22").Value = 1223.6316
$ws.Range("L22").Value = 2882.75
$ws.Range("M22").Value = -928.6315999999999
$ws.Range("N22").Value = -3472.75
$ws.Range("H27").Value = 1715.2222
$ws.Range("I27").Value = 1223.6316
$ws.Range("J27").Value = 2882.75
$ws.Range("K27").Value = 1223.6316
$ws.Range("L27").Value = 2882.75
$ws.Range("M27").Value = -1116.6316
$ws.Range("N27").Value = -3096.75
$ws.Range("H40").Value = 2204.818
$ws.Range("I40").Value = 2175.625
$ws.Range("J40").Value = 2282.6667
$ws.Range("K40").Value = 2175.625
$ws.Range("L40").Value = 2282.6667
$ws.Range("M40").Value = -2039.625
$ws.Range("N40").Value = -2554.6667
$ws.Range("H46").Value = 1417
$ws.Range("J46").Value = 1332.3334
$ws.Range("L46").Value = 1332.3334
$ws.Range("N46").Value = -1708.3334
$ws.Range("H48").Value = 3041
$ws.Range("I48").Value = 3041
$ws.Range("K48").Value = 3041
$ws.Range("M48").Value = -2380
$ws.Range("H50").Value = 54000
$ws.Range("J50").Value = 54000
$ws.Range("L50").Value = 54000
$ws.Range("N50").Value = -55274
$ws.Range("H55").Value = 487.44446
$ws.Range("I55").Value = 406.33334
$ws.Range("J55").Value = 649.6667
$ws.Range("K55").Value = 406.33334
$ws.Range("L55").Value = 649.6667
$ws.Range("M55").Value = -233.33334
$ws.Range("N55").Value = -995.6667
$ws.Range("H122").Value = 2850.125
$ws.Range("I122").Value = 2280.2
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 6840.599999999999
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -4390.599999999999
$ws.Range("N122").Value = -16300
$ws.Range("H126").Value = 2943.6667
$ws.Range("J126").Value = 2482.1667
$ws.Range("L126").Value = 7446.500100000001
$ws.Range("N126").Value = -12386.5001
$ws.Range("H132").Value = 1953.3334
$ws.Range("I132").Value = 1953.5883
$ws.Range("J132").Value = 1952.25
$ws.Range("K132").Value = 5860.7649
$ws.Range("L132").Value = 5856.75
$ws.Range("M132").Value = -3330.7649
$ws.Range("N132").Value = -10916.75
$ws.Range("H136").Value = 2876.5557
$ws.Range("I136").Value = 906.6667
$ws.Range("K136").Value = 2720.0001
$ws.Range("M136").Value = -170.0001000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 21494.25
$ws.Range("I74").Value = 18659
$ws.Range("K74").Value = 18659
$ws.Range("M74").Value = -17723
$ws.Range("H77").Value = 21494.25
$ws.Range("I77").Value = 18659
$ws.Range("K77").Value = 55977
$ws.Range("M77").Value = -51297
$ws.Range("H81").Value = 1279.8334
$ws.Range("I81").Value = 1279.8334
$ws.Range("K81").Value = 2559.6668
$ws.Range("M81").Value = -1498.6668
$ws.Range("H84").Value = 1279.8334
$ws.Range("I84").Value = 1279.8334
$ws.Range("K84").Value = 12798.334
$ws.Range("M84").Value = -7494.333999999999
$ws.Range("H122").Value = 11637.8
$ws.Range("I122").Value = 11637.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 34913.39999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -32463.39999999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2737.24
$ws.Range("I132").Value = 2525.3157
$ws.Range("J132").Value = 3408.3333
$ws.Range("K132").Value = 7575.9471
$ws.Range("L132").Value = 10224.9999
$ws.Range("M132").Value = -5045.9471
$ws.Range("N132").Value = -15284.9999
$ws.Range("H136").Value = 679.4
$ws.Range("I136").Value = 603.2963
$ws.Range("K136").Value = 1809.8889
$ws.Range("M136").Value = 740.1111000000001
$ws.Range("H138").Value = 67500
$ws.Range("J138").Value = 67500
$ws.Range("L138").Value = 67500
$ws.Range("N138").Value = -77780
